$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsValid = $wb.Worksheets.Item(1)    # was "Pythoncode"
$wsValid.Name = "ValidPythonCode"

$wsInvalid = $wb.Worksheets.Item(2)  # was "Sheet2"
$wsInvalid.Name = "InvalidPythonCode"

# --- ValidPythonCode: drop the stray 3rd row -----------------------------
$wsValid.Range("A3:B3").ClearContents()
[void]$wsValid.Range("B12").Select()

# --- InvalidPythonCode: fix up the error-message row + add new sample ---
$wsInvalid.Range("A2").Value = "Python programming 123 "
$wsInvalid.Range("B2").Value = "NameError: name 'Python' is not defined on line 1"

# match the bestFit column widths already used on the other sheet
$wsInvalid.Columns.Item(1).ColumnWidth = 11.666666666666666
$wsInvalid.Columns.Item(2).ColumnWidth = 42.166666666666664

[void]$wsInvalid.Range("B9").Select()
